# Updated symbol list on Sat Feb 11 03:30:03 UTC 2023 with GitHub Actions
# Applies the refreshed Price (column D) and Volume(1h) (column E) figures
# for the cryptos worksheet, preserving each cell's existing Text data type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "308.62" },
    @{ Cell = "E2"; Value = "0.33%" },
    @{ Cell = "D3"; Value = "41.47" },
    @{ Cell = "E3"; Value = "3.50%" },
    @{ Cell = "D4"; Value = "5.130" },
    @{ Cell = "E4"; Value = "2.46%" },
    @{ Cell = "D5"; Value = "0.07616" },
    @{ Cell = "E5"; Value = "-0.89%" },
    @{ Cell = "D6"; Value = "1.624" },
    @{ Cell = "E6"; Value = "-0.41%" },
    @{ Cell = "D7"; Value = "2.490" },
    @{ Cell = "E7"; Value = "-2.33%" },
    @{ Cell = "D8"; Value = "0.9045" },
    @{ Cell = "E8"; Value = "1.43%" },
    @{ Cell = "D9"; Value = "0.1079" },
    @{ Cell = "E9"; Value = "9.34%" },
    @{ Cell = "D10"; Value = "0.1766" },
    @{ Cell = "E10"; Value = "1.45%" },
    @{ Cell = "D11"; Value = "0.09138" },
    @{ Cell = "E11"; Value = "2.48%" },
    @{ Cell = "D12"; Value = "0.04263" },
    @{ Cell = "E12"; Value = "-2.55%" },
    @{ Cell = "D13"; Value = "0.1050" },
    @{ Cell = "E13"; Value = "-0.43%" },
    @{ Cell = "D14"; Value = "0.001257" },
    @{ Cell = "E14"; Value = "-1.14%" },
    @{ Cell = "D15"; Value = "0.005845" },
    @{ Cell = "E15"; Value = "-0.25%" },
    @{ Cell = "E16"; Value = "0.20%" },
    @{ Cell = "D17"; Value = "4.250" },
    @{ Cell = "E17"; Value = "0.10%" },
    @{ Cell = "D19"; Value = "6.532" },
    @{ Cell = "E19"; Value = "-7.12%" },
    @{ Cell = "D20"; Value = "0.1364" },
    @{ Cell = "E20"; Value = "0.99%" },
    @{ Cell = "D21"; Value = "0.2679" },
    @{ Cell = "E21"; Value = "-11.59%" },
    @{ Cell = "D22"; Value = "0.04192" },
    @{ Cell = "E22"; Value = "-0.85%" },
    @{ Cell = "D23"; Value = "0.001219" },
    @{ Cell = "E23"; Value = "1.60%" },
    @{ Cell = "D24"; Value = "0.004094" },
    @{ Cell = "E24"; Value = "0.61%" },
    @{ Cell = "D25"; Value = "0.0001299" },
    @{ Cell = "E25"; Value = "6.41%" },
    @{ Cell = "D38"; Value = "0.02411" },
    @{ Cell = "E38"; Value = "2.24%" },
    @{ Cell = "D39"; Value = "0.05194" },
    @{ Cell = "E39"; Value = "0.50%" },
    @{ Cell = "D40"; Value = "0.007773" },
    @{ Cell = "E40"; Value = "-1.98%" },
    @{ Cell = "D41"; Value = "0.1300" },
    @{ Cell = "E41"; Value = "-1.66%" },
    @{ Cell = "D42"; Value = "0.006947" },
    @{ Cell = "E42"; Value = "6.09%" },
    @{ Cell = "E43"; Value = "-5.82%" },
    @{ Cell = "D44"; Value = "0.008057" },
    @{ Cell = "E44"; Value = "5.64%" },
    @{ Cell = "D45"; Value = "0.3061" },
    @{ Cell = "E45"; Value = "-7.64%" },
    @{ Cell = "D46"; Value = "0.00006735" },
    @{ Cell = "E46"; Value = "1.45%" },
    @{ Cell = "E47"; Value = "-0.13%" },
    @{ Cell = "E48"; Value = "-12.00%" },
    @{ Cell = "D49"; Value = "0.009114" },
    @{ Cell = "E49"; Value = "190.88%" },
    @{ Cell = "E50"; Value = "-0.13%" },
    @{ Cell = "E51"; Value = "-0.13%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Keep the cell stored as text (matches the source data's inline-string
    # cells) instead of letting Excel auto-convert the numeric-looking /
    # percent-looking string into a Number cell.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

Write-Output "Updated $($updates.Count) cells"
